$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26: ID 24, Comment "I like this garden", Time "04-10-2023"
$ws.Range("A26").Value = 24
$ws.Range("B26").Value = "I like this garden"
$cellC26 = $ws.Range("C26")
$cellC26.Formula = "=""04-10-2023"""
$cellC26.Copy()
$cellC26.PasteSpecial(-4163)

# Row 27: ID 25, Comment "test id should be 25", Time "04-10-2023"
$ws.Range("A27").Value = 25
$ws.Range("B27").Value = "test id should be 25"
$cellC27 = $ws.Range("C27")
$cellC27.Formula = "=""04-10-2023"""
$cellC27.Copy()
$cellC27.PasteSpecial(-4163)

$excel.CutCopyMode = 0
